$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.043.73"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.465.29"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D5").Value = "'582.38"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'173.73"
$ws.Range("E6").Value = "  +2.96%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.511"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "2.913.01"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "66.921.77"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "2.468.28"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").Value = "'7.43"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").Value = "'347.94"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").Value = "'4.01"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D23").Value = "'69.27"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").Value = "'4.18"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").Value = "2.592.45"
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "0.0₃0895"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").Value = "'497.11"
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +2.53%  "
$ws.Range("D36").Value = "'160.55"
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'18.12"
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'142.27"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").Value = "'0.507"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").Value = "'0.580"
$ws.Range("E51").Value = "  +0.03%  "
